# Add a new match row (row 19) to the Gibraltar National League sheet,
# mirroring the layout/formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (styles/borders/number formats) of the last existing
# data row (18) onto the new row (19) before writing values into it.
$ws.Range("A18:V18").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A19").Value2 = 18
$ws.Range("B19").Value2 = "gibraltar"
$ws.Range("C19").Value2 = "national-league"
$ws.Range("D19").Value2 = "2023-2024"
$ws.Range("E19").Value2 = 45224.875
$ws.Range("F19").Value2 = "Magpies"
$ws.Range("G19").Value2 = 0
$ws.Range("H19").Value2 = "Lincoln Red Imps"
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 2.41
$ws.Range("K19").Value2 = "25/10/2023 11:34"
$ws.Range("L19").Value2 = 2.88
$ws.Range("M19").Value2 = "25/10/2023 19:51"
$ws.Range("N19").Value2 = 3.24
$ws.Range("O19").Value2 = "25/10/2023 11:34"
$ws.Range("P19").Value2 = 3.44
$ws.Range("Q19").Value2 = "25/10/2023 19:51"
$ws.Range("R19").Value2 = 2.58
$ws.Range("S19").Value2 = "25/10/2023 11:34"
$ws.Range("T19").Value2 = 2.14
$ws.Range("U19").Value2 = "25/10/2023 19:51"
$ws.Range("V19").Value2 = "https://www.betexplorer.com/football/gibraltar/national-league/magpies-lincoln-red-imps/APhWpNuB/"
